$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 43 with the new FAQ entry (Anahtar Kelime, Senaryo, Açıklama,
# Çözüm, Sorumlu, Görsel) - this introduces column F data for the first time.
$ws.Cells.Item(43, 1).Value = "teslim"
$ws.Cells.Item(43, 2).Value = "Aracı kime teslim edeceğim"
$ws.Cells.Item(43, 3).Value = "Aracı teslim edeceğin kişi teslimat sekmesi altında yer alır. Bu alanı ihaleye girip ödemesini yapan ve mail adresi sistemde kayıtlı olan kişi doldurur."
$ws.Cells.Item(43, 4).Value = "Salesforce da plakayı yazıp ara. Teslimat aşamasında olan kaydı seç. Teslimat sekmesini seç ve sağ alttaki kişi ve TC kimlik numarasını kontrol et."
$ws.Cells.Item(43, 5).Value = "Product Manager"
$ws.Cells.Item(43, 6).Value = "teslim edilecek kişi.JPG"

# Scroll the view down and move the active selection, matching the saved
# workbook view state (topLeftCell="A22", selection at G26).
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("G26").Select()
